$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

$ws.Range("B23").Value = "Micro: <5, `nVery small: `n<10 Agriculture `n<20 Mining, Manuf., Elec., ConsTurnover, Retail, Wholesale, Catering, Transport, Finan., Others"
$ws.Range("C23").Value = "Micro: Rm<0.10 All sectors`nVery small: `nRm<0.5 Agriculture Construction Finance `n<2 Mining, Manuf.,`n<1.90 Elec., Catering `n<0.6 Retail, Wholesale, TranAnnual sales, Others"
$ws.Range("D23").Value = "Micro: Rm<0.20 All sectors`nVery small: `nRm<0.5 Agriculture `n<4 Mining, Retail `n<52 Manuf., `n<5.1 Elec., `n<3 ConsTurnover, Transport Finance`n <6 Wholesale, `n<5.1 Catering,`n<1 Others"
$ws.Range("C24").Value = "<3 Rm Agriculture, `n<6 Mining, `n<5 Manuf., Elec., Wholesale `n<1 ConsTurnover and Catering `n<3 Retail, TranAnnual sales, Finance, OtherAnnual sales"
$ws.Range("D24").Value = "<3 Rm Agriculture, `n<10 Mining, `n<13 Manuf. TranAnnual sales and Elec., `n<6 ConsTurnover and catering `n<19 Retail, `n<32 Wholesale"
$ws.Range("B25").Value = "<100 Agri. `n<200 Others"
$ws.Range("C25").Value = "<5 Rm Agriculture Finance, `n<23 Mining, `n<19 Manuf. and Elec., `n<5 Const `n<6 Retail, Transport and Others `n<10 Wholesale, `n<3 Catering"
$ws.Range("D25").Value = "<5 Rm Agriculture,`n<39 Mining, `n<51 Manuf. and Elec., `n<26 Const, Transport `n<39 Retail, `n<64 Wholesale, `n<13 Catering"
$ws.Range("B26").Value = ">=100 Agri. `n>=200 Others"
$ws.Range("C26").Value = ">=5 Rm Agriculture Finance, `n>=23 Mining,`n >=19 Manuf. and Elec., `n>=5 Const >=6 Retail, Transport and Others `n>=10 Wholesale, `n>=3 Catering"
$ws.Range("D26").Value = ">=5 Rm Agriculture, `n>=39 Mining, `n>=51 Manuf. and Elec.,`n>=26 Const, Transport `n>=39 Retail, `n>=64 Wholesale, `n>=13 Catering"
